$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Athletic Club
$ws.Range("B3").Value = 27
$ws.Range("C3").Value = 28
$ws.Range("D3").Value = 49.5
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 275
$ws.Range("G3").Value = 2250
$ws.Range("H3").Value = 25
$ws.Range("I3").Value = 28
$ws.Range("J3").Value = 17
$ws.Range("K3").Value = 45
$ws.Range("L3").Value = 23
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = 58
$ws.Range("Q3").Value = 1.12
$ws.Range("R3").Value = 0.68
$ws.Range("S3").Value = 1.8
$ws.Range("U3").Value = 1.6

# Row 7 - Elche
$ws.Range("B7").Value = 32
$ws.Range("D7").Value = 58.4
$ws.Range("E7").Value = 25
$ws.Range("F7").Value = 275
$ws.Range("G7").Value = 2250
$ws.Range("H7").Value = 25
$ws.Range("I7").Value = 31
$ws.Range("K7").Value = 54
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 52
$ws.Range("Q7").Value = 1.24
$ws.Range("R7").Value = 0.92
$ws.Range("S7").Value = 2.16
$ws.Range("T7").Value = 1.2
$ws.Range("U7").Value = 2.12
